# Adds "Romania" and "Slovakia" market test-data sheets between "Belgium"
# and "Turkey", inserts two new rows (MX-BBX / MX-DPBX) into every sheet's
# accessories list, and makes "Turkey" the active tab - mirroring the
# "Added RomaniaFC,SlovakiaFC test data" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the two new country sheets by copying "Belgium" (same layout,
#    styles, column widths, merged cells) and renaming/placing them.
# ---------------------------------------------------------------------
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Copy($null, $belgium)
$romania = $wb.Worksheets.Item("Belgium (2)")
$romania.Name = "Romania"

$romania.Copy($null, $romania)
$slovakia = $wb.Worksheets.Item("Romania (2)")
$slovakia.Name = "Slovakia"

# ---------------------------------------------------------------------
# 2. Fill in the market name / model-number headers for the new sheets.
#    Order matters here: it controls the order new shared strings are
#    interned in, which must match "Romania Market", "MX-DPBX",
#    "NGC-4307/T3536/T3543", "NGC-4306/T3562/T3575", "Slovakia Market",
#    "MX-BBX".
# ---------------------------------------------------------------------
$romania.Range("B2").Value = "Romania Market"

$uk = $wb.Worksheets.Item("UK")
$belgiumWs = $wb.Worksheets.Item("Belgium")
$turkey = $wb.Worksheets.Item("Turkey")

$allSheets = @($uk, $belgiumWs, $romania, $slovakia, $turkey)

# ---------------------------------------------------------------------
# 3. Insert the two new accessory rows ("MX-BBX" row then "MX-DPBX" row)
#    above the trailing "FC1D2-Unmonitored" / "Wg" / "Accessories" rows
#    on every sheet, copying the formatting down from the row above so
#    the cell style ("s=3") matches.
# ---------------------------------------------------------------------
foreach ($ws in $allSheets) {
    $ws.Rows.Item(12).Insert()
    $ws.Rows.Item(12).Insert()
    $ws.Range("A11").Copy()
    $ws.Range("A12").PasteSpecial(-4122)
    $ws.Range("A13").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Write the "MX-DPBX" row text first (row 13) on every sheet - this is the
# first brand-new string Excel will intern after "Romania Market".
foreach ($ws in $allSheets) {
    $ws.Range("A13").Value = "MX-DPBX"
}

# Romania / Slovakia model numbers + Slovakia market name, in order.
$romania.Range("B4").Value = "NGC-4307/T3536/T3543"
$slovakia.Range("B4").Value = "NGC-4306/T3562/T3575"
$slovakia.Range("B2").Value = "Slovakia Market"

# Finally the "MX-BBX" row text (row 12) on every sheet - last new string.
foreach ($ws in $allSheets) {
    $ws.Range("A12").Value = "MX-BBX"
}

# ---------------------------------------------------------------------
# 4. Selection + active tab bookkeeping to match the target workbook:
#    every sheet's cursor rests on the freshly-added "A12" row, and the
#    last tab ("Turkey") becomes the active one.
# ---------------------------------------------------------------------
foreach ($ws in $allSheets) {
    $ws.Range("A12").Select()
}

$turkey.Activate()
